# Automatic update of files.
# Applies the recorded edits to the "Artfynd" sheet:
#  - Round the Ost (Q) / Nord (R) coordinate values that carried long
#    floating point tails down to whole metres.
#  - Rows 3/4, 7/8 and 9/10/11 had their observation details re-matched to
#    the correct location/time/taxon - update each cell to the corrected
#    value.
#  - A few Starttid/Sluttid ("00:00") cells turned out to be blank, and a
#    couple of Enhet ("plantor/tuvor") / Antal values moved to a different
#    row along with the rest of their record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: round Ost/Nord -------------------------------------------------
$ws.Range("Q2").Value = 489818
$ws.Range("R2").Value = 6949032

# --- Rows 3 & 4: the two observations had swapped Id/time/coordinates -----
$ws.Range("A3").Value = 111792337
$ws.Range("Q3").Value = 489764
$ws.Range("R3").Value = 6949092
$ws.Range("Z3").Value = "19:22"
$ws.Range("AB3").Value = "19:22"

$ws.Range("A4").Value = 111790625
$ws.Range("Q4").Value = 489825
$ws.Range("R4").Value = 6949021
$ws.Range("Z4").Value = "18:29"
$ws.Range("AB4").Value = "18:29"

# --- Row 5: round Ost/Nord, clear the (blank) time cells -------------------
$ws.Range("Q5").Value = 490133
$ws.Range("R5").Value = 6948774
$ws.Range("Z5").ClearContents()
$ws.Range("AB5").ClearContents()

# --- Row 6: round Ost/Nord, clear the (blank) time cells -------------------
$ws.Range("Q6").Value = 490124
$ws.Range("R6").Value = 6948875
$ws.Range("Z6").ClearContents()
$ws.Range("AB6").ClearContents()

# --- Rows 7 & 8: the two observations had swapped Id/location/time/enhet --
$ws.Range("A7").Value = 112212882
$ws.Range("J7").Value = "plantor/tuvor"
$ws.Range("P7").Value = "Kälen (Kälen), Jmt"
$ws.Range("Q7").Value = 490109
$ws.Range("R7").Value = 6948768
$ws.Range("Z7").Value = "12:39"
$ws.Range("AB7").Value = "12:39"

$ws.Range("A8").Value = 112212105
$ws.Range("J8").ClearContents()
$ws.Range("P8").Value = "Nordvallen (Nordvallen), Jmt"
$ws.Range("Q8").Value = 490018
$ws.Range("R8").Value = 6948882
$ws.Range("Z8").Value = "11:58"
$ws.Range("AB8").Value = "11:58"

# --- Rows 9, 10 & 11: taxon/location data rotated across the three rows ---
$ws.Range("A9").Value = 112212902
$ws.Range("B9").Value = 78578
$ws.Range("E9").Value = 6458
$ws.Range("F9").Value = "Lunglav"
$ws.Range("G9").Value = "Lobaria pulmonaria"
$ws.Range("H9").Value = "(L.) Hoffm."
$ws.Range("P9").Value = "Kälen (Kälen), Jmt"
$ws.Range("Q9").Value = 490134
$ws.Range("R9").Value = 6948772

$ws.Range("A10").Value = 112213279
$ws.Range("B10").Value = 89405
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 1202
$ws.Range("F10").Value = "Ullticka"
$ws.Range("G10").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H10").Value = "(P.Karst.) Fiasson & Niemelä"
# Antal (I10) keeps its text formatting but the value moved to I11 - leave
# this cell blank (text-typed, not a plain cleared cell).
$ws.Range("I10").NumberFormat = "@"
$ws.Range("I10").Value = ""
$ws.Range("J10").ClearContents()
$ws.Range("P10").Value = "Nordvallen (Nordvallen), Jmt"
$ws.Range("Q10").Value = 490080
$ws.Range("R10").Value = 6948907

$ws.Range("A11").Value = 112212836
$ws.Range("B11").Value = 96348
$ws.Range("D11").Value = "VU"
$ws.Range("E11").Value = 220787
$ws.Range("F11").Value = "Knärot"
$ws.Range("G11").Value = "Goodyera repens"
$ws.Range("H11").Value = "(L.) R. Br."
# Antal (I11) is text "25", not the number 25 - force text formatting so it
# round-trips as a string like the source record did.
$ws.Range("I11").NumberFormat = "@"
$ws.Range("I11").Value = "25"
$ws.Range("J11").Value = "plantor/tuvor"
$ws.Range("P11").Value = "Stugunäset (Stugunäset), Jmt"
$ws.Range("Q11").Value = 490078
$ws.Range("R11").Value = 6948752
